# feat: improve tool look
#
# The "50 controls" performance measurement (row 7) is updated from 3689s
# to 4038s, and the active selection is moved off the stale D11 reference
# back onto the data table (A2).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the measured time for n=50 controls.
$ws.Range("B7").Value = 4038

# Move the selection back onto the data range instead of the old D11 cell.
$ws.Range("A2").Select()
